$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of user stories for Sprint 4 / NOT STARTED
$ws.Range("A11").Value = "SPRINT 4"
$ws.Range("B11").Value = "I want to get avatar for each housemate from the phone contact"
$ws.Range("C11").Value = "NOT STARTED"

$ws.Range("A12").Value = "SPRINT 4"
$ws.Range("B12").Value = "I want to dispay room avatar as a combined image from roommates' avatars"
$ws.Range("C12").Value = "NOT STARTED"

# Match the "Accent3" cell style used by the Sprint/Status columns in the row above
$ws.Range("A11").Style = "Accent3"
$ws.Range("C11").Style = "Accent3"
$ws.Range("A12").Style = "Accent3"
$ws.Range("C12").Style = "Accent3"

# Resize the table (ListObject) to include the two new rows
$table = $ws.ListObjects.Item("Table2")
$table.Resize($ws.Range("A1:C12"))

# Widen column B to fit the new, longer text (matches Excel's computed best-fit width)
$ws.Columns.Item(2).ColumnWidth = 68.25

# Update the selected cell
$ws.Range("B15").Select()
